$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit re-orders the species records currently sitting in rows 19-26
# (same locality/date/observer batch) and bumps a handful of
# "Taxonsorteringsordning" (col B) values by +14 along the way. The
# "Publik kommentar" note in AC travels with its parent record (the
# Picoides tridactylus / Tretaig hackspett row), moving from row 20 to
# row 22.

# Row 19 <- old row 26 (Skrovellav / Lobaria scrobiculata), B bumped
$ws.Range("A19").Value = 112435800
$ws.Range("B19").Value = 78714
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 2081
$ws.Range("F19").Value = "Skrovellav"
$ws.Range("G19").Value = "Lobaria scrobiculata"
$ws.Range("H19").Value = "(Scop.) DC."
$ws.Range("Q19").Value = 428438
$ws.Range("R19").Value = 6967379

# Row 20 <- old row 19 (Bardlav / Nephroma parile), B bumped
$ws.Range("A20").Value = 112435711
$ws.Range("B20").Value = 78746
$ws.Range("D20").Value = "LC"
$ws.Range("E20").Value = 6463
$ws.Range("F20").Value = "Bårdlav"
$ws.Range("G20").Value = "Nephroma parile"
$ws.Range("H20").Value = "(Ach.) Ach."
$ws.Range("Q20").Value = 428423
$ws.Range("R20").Value = 6967395
$ws.Range("AC20").ClearContents()

# Row 21 <- old row 22 (Stuplav / Nephroma bellum), B bumped
$ws.Range("A21").Value = 112435836
$ws.Range("B21").Value = 78740
$ws.Range("D21").Value = "LC"
$ws.Range("E21").Value = 6462
$ws.Range("F21").Value = "Stuplav"
$ws.Range("G21").Value = "Nephroma bellum"
$ws.Range("H21").Value = "(Spreng.) Tuck."
$ws.Range("Q21").Value = 428438
$ws.Range("R21").Value = 6967379

# Row 22 <- old row 20 (Tretaig hackspett / Picoides tridactylus), B unchanged,
# and the public comment follows this record into row 22.
$ws.Range("A22").Value = 112439263
$ws.Range("B22").Value = 56430
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 100109
$ws.Range("F22").Value = "Tretåig hackspett"
$ws.Range("G22").Value = "Picoides tridactylus"
$ws.Range("H22").Value = "(Linnaeus, 1758)"
$ws.Range("Q22").Value = 428455
$ws.Range("R22").Value = 6967575
$ws.Range("AC22").Value = "Gamla ringhack på tall."

# Row 23 <- old row 24 (Spadskinn / Stereopsis vitellina), B bumped
$ws.Range("A23").Value = 112438686
$ws.Range("B23").Value = 90795
$ws.Range("D23").Value = "VU"
$ws.Range("E23").Value = 6055
$ws.Range("F23").Value = "Spadskinn"
$ws.Range("G23").Value = "Stereopsis vitellina"
$ws.Range("H23").Value = "(S.Lundell) D.A.Reid"
$ws.Range("Q23").Value = 428467
$ws.Range("R23").Value = 6967562
$ws.Range("S23").Value = 10

# Row 24 <- old row 21 (Svartvit taggsvamp / Phellodon connatus), B bumped
$ws.Range("A24").Value = 112439218
$ws.Range("B24").Value = 90857
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 5448
$ws.Range("F24").Value = "Svartvit taggsvamp"
$ws.Range("G24").Value = "Phellodon connatus"
$ws.Range("H24").Value = "(Schultz) nom.prov"
$ws.Range("Q24").Value = 428439
$ws.Range("R24").Value = 6967600
$ws.Range("S24").Value = 20

# Row 25 <- old row 23 (Skrovellav / Lobaria scrobiculata), B bumped
$ws.Range("A25").Value = 112435620
$ws.Range("B25").Value = 78714
$ws.Range("E25").Value = 2081
$ws.Range("F25").Value = "Skrovellav"
$ws.Range("G25").Value = "Lobaria scrobiculata"
$ws.Range("H25").Value = "(Scop.) DC."
$ws.Range("Q25").Value = 428423
$ws.Range("R25").Value = 6967395

# Row 26 <- old row 25 (Mork kolflarnlav / Carbonicola myrmecina), B bumped
$ws.Range("A26").Value = 112436511
$ws.Range("B26").Value = 77403
$ws.Range("E26").Value = 228912
$ws.Range("F26").Value = "Mörk kolflarnlav"
$ws.Range("G26").Value = "Carbonicola myrmecina"
$ws.Range("H26").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q26").Value = 428502
$ws.Range("R26").Value = 6967309
